$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.833.53'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.874.82'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4587'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3869'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07867'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9832'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.70'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.891.92'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.981'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.644'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06953'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009958'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '28.845.08'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.237'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.982'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.922'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09328'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9002'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.252'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.315'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.184'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05753'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02063'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.002'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.653'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5633'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1760'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.610'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.254'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5335'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07038'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.838'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '112.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.507'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.061'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.52%  '
